$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.054260969161987
$ws.Range("B1").Value = 2.241319417953491
$ws.Range("C1").Value = 2.304770946502686
$ws.Range("D1").Value = 2.80654764175415
$ws.Range("E1").Value = 3.364995718002319
